# Challenge timesheet update ("ch12 - multiple check box completed"):
# Fill in the "Name" column (D) for rows 8-13 with the challenge titles,
# and mark rows 8-11 (Challenge-07..Challenge-10) as Completed since
# those challenges are now finished. Rows 12-13 keep their Pending status.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8  (Challenge-07): name + mark Completed
$ws.Range("D8").Value = "Array Works 1"
$ws.Range("E8").Value = "Completed"

# Row 9  (Challenge-08): name + mark Completed
$ws.Range("D9").Value = "Canvas"
$ws.Range("E9").Value = "Completed"

# Row 10 (Challenge-09): name + mark Completed
$ws.Range("D10").Value = "14 dev tools"
$ws.Range("E10").Value = "Completed"

# Row 11 (Challenge-10): name + mark Completed
$ws.Range("D11").Value = "multiple check box"
$ws.Range("E11").Value = "Completed"

# Row 12 (Challenge-11): name only, remains Pending
$ws.Range("D12").Value = "Cuztom html5 video player"

# Row 13 (Challenge-12): name only, remains Pending
$ws.Range("D13").Value = "Key Sequence detection"

# Update the last active selection to match the authored state
$ws.Range("K12").Select()
